# "Generate Report for Archive"
# The localization-status report is regenerated: the in-flight status text
# changes from "Ready for handoff" to "In Translation" on every sheet that
# shows it (Overview!E2/F2, zh-cn!C2, de-de!C2). Because the new status
# string is shorter than the old one, the status columns that were sized to
# fit the old text shrink to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status cells: "Ready for handoff" -> "In Translation"
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Shrink the (now narrower) status columns to fit the new text.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
